$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(13, 'Eastern Michigan', 'Akron', 0, 3, 1, 0),
    @(13, 'Northern Illinois', 'Western Michigan', 0, 24, 1, 0),
    @(13, 'Bowling Green', 'Toledo', 0, 1, 0, 1),
    @(13, 'Miami (OH)', 'Buffalo', 0, 13, 1, 0),
    @(13, 'Ohio', 'Central Michigan', 0, 14, 1, 0),
    @(13, 'Pittsburgh', 'Boston College', 0, 8, 1, 0),
    @(13, 'Texas-San Antonio', 'South Florida', 0, 28, 1, 0),
    @(13, 'Washington State', 'Colorado', 0, 42, 1, 0),
    @(13, 'Alabama', 'Chattanooga', 0, 56, 1, 0),
    @(13, 'Alabama-Birmingham', 'Temple', 0, 10, 1, 0),
    @(13, 'James Madison', 'Appalachian State', 0, 3, 0, 1),
    @(13, 'Arizona', 'Utah', 0, 24, 1, 0),
    @(13, 'Arkansas', 'Florida International', 0, 24, 1, 0),
    @(13, 'Arkansas State', 'Texas State', 0, 46, 1, 0),
    @(13, 'Army', 'Coastal Carolina', 0, 7, 1, 0),
    @(13, 'Ball State', 'Kent State', 0, 31, 1, 0),
    @(13, 'Utah State', 'Boise State', 0, 35, 0, 1),
    @(13, 'Stanford', 'California', 0, 12, 0, 1),
    @(13, 'Clemson', 'North Carolina', 0, 11, 1, 0),
    @(13, 'Colorado State', 'Nevada', 0, 10, 1, 0),
    @(13, 'Connecticut', 'Sacred Heart', 0, 28, 1, 0),
    @(13, 'Florida State', 'North Alabama', 0, 45, 1, 0),
    @(13, 'Tennessee', 'Georgia', 0, 28, 0, 1),
    @(13, 'Georgia Tech', 'Syracuse', 0, 9, 1, 0),
    @(13, 'Iowa', 'Illinois', 0, 2, 1, 0),
    @(13, 'Jacksonville State', 'Louisiana Tech', 0, 39, 1, 0),
    @(13, 'Kansas', 'Kansas State', 0, 4, 0, 1),
    @(13, 'Liberty', 'Massachusetts', 0, 24, 1, 0),
    @(13, 'Louisiana State', 'Georgia State', 0, 42, 1, 0),
    @(13, 'Miami (FL)', 'Louisville', 0, 7, 0, 1),
    @(13, 'Maryland', 'Michigan', 0, 7, 0, 1),
    @(13, 'Indiana', 'Michigan State', 0, 3, 0, 1),
    @(13, 'Middle Tennessee State', 'Texas-El Paso', 0, 4, 1, 0),
    @(13, 'Mississippi', 'Louisiana-Monroe', 0, 32, 1, 0),
    @(13, 'Mississippi State', 'Southern Mississippi', 0, 21, 1, 0),
    @(13, 'Missouri', 'Florida', 0, 2, 1, 0),
    @(13, 'Navy', 'East Carolina', 0, 10, 1, 0),
    @(13, 'Air Force', 'Nevada-Las Vegas', 0, 4, 0, 1),
    @(13, 'Fresno State', 'New Mexico', 0, 8, 0, 1),
    @(13, 'Auburn', 'New Mexico State', 0, 21, 0, 1),
    @(13, 'Virginia Tech', 'North Carolina State', 0, 7, 0, 1),
    @(13, 'Tulsa', 'North Texas', 0, 7, 0, 1),
    @(13, 'Northwestern', 'Purdue', 0, 8, 1, 0),
    @(13, 'Notre Dame', 'Wake Forest', 0, 38, 1, 0),
    @(13, 'Ohio State', 'Minnesota', 0, 34, 1, 0),
    @(13, 'Brigham Young', 'Oklahoma', 0, 7, 0, 1),
    @(13, 'Houston', 'Oklahoma State', 0, 13, 0, 1),
    @(13, 'Georgia Southern', 'Old Dominion', 0, 3, 0, 1),
    @(13, 'Arizona State', 'Oregon', 0, 36, 0, 1),
    @(13, 'Penn State', 'Rutgers', 0, 21, 1, 0),
    @(13, 'Charlotte', 'Rice', 0, 21, 0, 1),
    @(13, 'San Jose State', 'San Diego State', 0, 11, 1, 0),
    @(13, 'South Alabama', 'Marshall', 0, 28, 1, 0),
    @(13, 'South Carolina', 'Kentucky', 0, 3, 1, 0),
    @(13, 'Memphis', 'Southern Methodist', 0, 4, 0, 1),
    @(13, 'Iowa State', 'Texas', 0, 10, 0, 1),
    @(13, 'Texas A&M', 'Abilene Christian', 0, 28, 1, 0),
    @(13, 'Texas Christian', 'Baylor', 0, 25, 1, 0),
    @(13, 'Texas Tech', 'Central Florida', 0, 1, 1, 0),
    @(13, 'Troy', 'Louisiana', 0, 7, 1, 0),
    @(13, 'Florida Atlantic', 'Tulane', 0, 16, 0, 1),
    @(13, 'Southern California', 'UCLA', 0, 18, 0, 1),
    @(13, 'Virginia', 'Duke', 0, 3, 1, 0),
    @(13, 'Oregon State', 'Washington', 0, 2, 0, 1),
    @(13, 'West Virginia', 'Cincinnati', 0, 21, 1, 0),
    @(13, 'Western Kentucky', 'Sam Houston', 0, 5, 1, 0),
    @(13, 'Wisconsin', 'Nebraska', 0, 7, 1, 0),
    @(13, 'Wyoming', 'Hawaii', 0, 33, 1, 0)
)

$startRow = 726
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $item = $data[$i]
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = ""
    $ws.Cells.Item($row, 5).Value = $item[3]
    $ws.Cells.Item($row, 6).Value = $item[4]
    $ws.Cells.Item($row, 7).Value = $item[5]
    $ws.Cells.Item($row, 8).Value = $item[6]
}
